$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "1" to "Abasha"
$ws.Name = "Abasha"

# Urban row (row 6): mark 2012, 2013, 2015, 2016, 2019 as confidential/unavailable
$ws.Range("D6").Value = "…"
$ws.Range("E6").Value = "…"
$ws.Range("G6").Value = "…"
$ws.Range("H6").Value = "…"
$ws.Range("K6").Value = "…"

# Rural row (row 7): mark 2012, 2013, 2019 as confidential/unavailable
$ws.Range("D7").Value = "…"
$ws.Range("E7").Value = "…"
$ws.Range("K7").Value = "…"

# Remove the blank spacer row between the data table and the note (old row 8),
# shifting the "Note" row up from row 9 to row 8
$ws.Rows(8).Delete()
